# Append the new "Consultar un producto en el inventario por estado" SQL
# block at the end of the document (right after the existing "LIMIT 1"; "
# paragraph, before the sectPr). We build the addition as a WordprocessingML
# fragment (wrapped in the pkg:package form that Range.WordOpenXML/InsertXML
# use) so the generated <w:p>/<w:r>/<w:proofErr> structure matches exactly
# what Word would produce when this text was typed in, including the
# spell-check proofErr wrappers around the dotted identifiers.

$d = $word.ActiveDocument

$insertXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>Consultar un producto en el inventario por estado</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.detalle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.estado</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.fecha_registro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.ingreso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.total</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.username</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">,  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.nombres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.nombre_producto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.inventario</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> i</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.usuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> u </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.proveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pro </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.idproveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.idProveedor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.productos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> p </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i.estado</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 'Ingresado';</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)
$null = $insertionPoint.InsertXML($insertXml)
